# Apply the "reference doc" style updates described by the golden-test diff:
#   - Title / TitleChar: drop the expanded/condensed character spacing and
#     minimum-kerning-size overrides from the run properties.
#   - Author / Date: re-base on Title (so they inherit its look), and carry
#     their own explicit run size (12pt) instead of relying on a direct
#     paragraph alignment override.

$d = $word.ActiveDocument

# --- Title ---------------------------------------------------------------
$title = $d.Styles("Title")
$title.Font.Spacing = 0
$title.Font.Kerning = 0

# --- TitleChar (the linked character style for Title) ---------------------
$titleChar = $d.Styles("TitleChar")
$titleChar.Font.Spacing = 0
$titleChar.Font.Kerning = 0

# --- Author ----------------------------------------------------------------
$author = $d.Styles("Author")
$author.BaseStyle = "Title"
$author.Font.Size = 12
$author.Font.SizeBi = 12

# --- Date --------------------------------------------------------------
$date = $d.Styles("Date")
$date.BaseStyle = "Title"
$date.Font.Size = 12
$date.Font.SizeBi = 12
